{"js": "// \"Se ha modificado el mensaje\": insert a new paragraph reading\n// \"Aprendizaje en Git \" immediately after the \"Prueba Git \u2013 modificado por\n// mateus\" paragraph (and before the pre-existing trailing empty paragraph).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that carries the \"modificado por\" text; fall back to\n// the very first paragraph if, for some reason, it can't be found.\nlet targetParagraph = paragraphs.items.find((p) =>\n  p.text.indexOf(\"modificado por\") !== -1\n);\nif (!targetParagraph) {\n  targetParagraph = paragraphs.items[0];\n}\n\ntargetParagraph.insertParagraph(\"Aprendizaje en Git \", \"After\");\nawait context.sync();\n", "ps1": "# \"Se ha modificado el mensaje\": insert a new paragraph reading\n# \"Aprendizaje en Git \" immediately after the \"Prueba Git - modificado por\n# mateus\" paragraph (and before the pre-existing trailing empty paragraph).\n$d = $word.ActiveDocument\n\n# Locate the paragraph that carries the \"modificado por\" text.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"modificado por\")\nif ($found) {\n    $targetPara = $searchRange.Paragraphs(1)\n} else {\n    $targetPara = $d.Paragraphs(1)\n}\n\n$tailRange = $targetPara.Range\n$tailRange.Collapse(0)          # wdCollapseEnd\n$tailRange.InsertParagraphAfter()\n$tailRange.InsertAfter(\"Aprendizaje en Git \")\n"}
